$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear column F (old rightmost column) since data is shifting left by one column
$ws.Range("F1:F3").Clear()

# A1 is a brand-new cell; copy the header formatting (border/bold/center) from B1 first
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122) | Out-Null

# Row 1 - headers (shifted left by one column, keep style)
$ws.Range("A1").Value = "QS_Phylonet5"
$ws.Range("B1").Value = "FNRATE_EXACT_ASTRAL"
$ws.Range("C1").Value = "TAXON"
$ws.Range("D1").Value = "MODEL_CONDITION"
$ws.Range("E1").Value = "GENE"

# Row 2
$ws.Range("A2").Value = 1508
$ws.Range("B2").Value = 0.125
$ws.Range("C2").Value = "11-texon"
$ws.Range("D2").Value = "estimated_5genes_weakILS"
$ws.Range("E2").Value = 4

# Row 3
$ws.Range("A3").Value = 1508
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = "11-texon"
$ws.Range("D3").Value = "estimated_5genes_weakILS"
$ws.Range("E3").Value = 15

# Old A2/A3 had style "1" (bordered) that must not carry forward to the new A2/A3 (now plain)
$ws.Range("A2:A3").ClearFormats()
